$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 140689.67
$ws.Range("J19").Value = 140689.67
$ws.Range("L19").Value = 140689.67
$ws.Range("N19").Value = -141039.67
# Row 131
$ws.Range("H131").Value = 457004.62
$ws.Range("I131").Value = 501905.6
$ws.Range("K131").Value = 1505716.8
$ws.Range("M131").Value = -1500676.8
# Row 132
$ws.Range("H132").Value = 50971.5
$ws.Range("I132").Value = 96668.336
$ws.Range("K132").Value = 290005.008
$ws.Range("M132").Value = -287475.008
# Row 135
$ws.Range("H135").Value = 2780.8235
$ws.Range("J135").Value = 1750
$ws.Range("L135").Value = 15750
$ws.Range("N135").Value = -20820
# Row 137
$ws.Range("H137").Value = 12581.105
$ws.Range("J137").Value = 31841.572
$ws.Range("L137").Value = 95524.716
$ws.Range("N137").Value = -100624.716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4558.161
$ws.Range("I2").Value = 3387.5217
$ws.Range("J2").Value = 7923.75
$ws.Range("K2").Value = 3387.5217
$ws.Range("L2").Value = 7923.75
$ws.Range("M2").Value = -3274.5217
$ws.Range("N2").Value = -8149.75
# Row 61
$ws.Range("H61").Value = 82580.34
$ws.Range("I61").Value = 2914.7856
$ws.Range("J61").Value = 156934.86
$ws.Range("K61").Value = 2914.7856
$ws.Range("L61").Value = 156934.86
$ws.Range("M61").Value = -2702.7856
$ws.Range("N61").Value = -157358.86
# Row 74
$ws.Range("H74").Value = 11146.313
$ws.Range("I74").Value = 1318.7906
$ws.Range("J74").Value = 63969.25
$ws.Range("K74").Value = 1318.7906
$ws.Range("L74").Value = 63969.25
$ws.Range("M74").Value = -444.7906
$ws.Range("N74").Value = -65717.25
# Row 77
$ws.Range("H77").Value = 11146.313
$ws.Range("I77").Value = 1318.7906
$ws.Range("J77").Value = 63969.25
$ws.Range("K77").Value = 6593.953
$ws.Range("L77").Value = 319846.25
$ws.Range("M77").Value = -2225.953
$ws.Range("N77").Value = -328582.25
# Row 88
$ws.Range("H88").Value = 2857.4285
$ws.Range("I88").Value = 2667.3333
$ws.Range("K88").Value = 2667.3333
$ws.Range("M88").Value = -2261.3333
# Row 91
$ws.Range("H91").Value = 2857.4285
$ws.Range("I91").Value = 2667.3333
$ws.Range("K91").Value = 2667.3333
$ws.Range("M91").Value = -1263.3333
# Row 97
$ws.Range("H97").Value = 812.6429000000001
$ws.Range("I97").Value = 812.6429000000001
$ws.Range("K97").Value = 812.6429000000001
$ws.Range("M97").Value = -316.6429000000001
# Row 116
$ws.Range("H116").Value = 4558.161
$ws.Range("I116").Value = 3387.5217
$ws.Range("J116").Value = 7923.75
$ws.Range("K116").Value = 3387.5217
$ws.Range("L116").Value = 7923.75
$ws.Range("M116").Value = -1093.5217
$ws.Range("N116").Value = -12511.75
# Row 122
$ws.Range("H122").Value = 1065366.2
$ws.Range("I122").Value = 1383302
$ws.Range("K122").Value = 4149906
$ws.Range("M122").Value = -4147456
# Row 132
$ws.Range("H132").Value = 2788878
$ws.Range("I132").Value = 2958.926
$ws.Range("K132").Value = 8876.778
$ws.Range("M132").Value = -6346.778
# Row 136
$ws.Range("H136").Value = 82580.34
$ws.Range("I136").Value = 2914.7856
$ws.Range("J136").Value = 156934.86
$ws.Range("K136").Value = 8744.356800000001
$ws.Range("L136").Value = 470804.58
$ws.Range("M136").Value = -6194.356800000001
$ws.Range("N136").Value = -475904.58

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4558.161
$ws.Range("I3").Value = 3387.5217
$ws.Range("J3").Value = 7923.75
$ws.Range("K3").Value = 3387.5217
$ws.Range("L3").Value = 7923.75
$ws.Range("M3").Value = -3273.5217
$ws.Range("N3").Value = -8151.75
# Row 80
$ws.Range("H80").Value = 1135.2778
$ws.Range("J80").Value = 1205.25
$ws.Range("L80").Value = 1205.25
$ws.Range("N80").Value = -3201.25
# Row 83
$ws.Range("H83").Value = 1135.2778
$ws.Range("J83").Value = 1205.25
$ws.Range("L83").Value = 6026.25
$ws.Range("N83").Value = -16010.25
# Row 99
$ws.Range("H99").Value = 10112
$ws.Range("I99").Value = 10989.92
$ws.Range("J99").Value = 4625
$ws.Range("K99").Value = 10989.92
$ws.Range("L99").Value = 4625
$ws.Range("M99").Value = -9491.92
$ws.Range("N99").Value = -7621
# Row 107
$ws.Range("H107").Value = 2381.45
$ws.Range("I107").Value = 2476.8125
$ws.Range("K107").Value = 2476.8125
$ws.Range("M107").Value = -556.8125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 11424.625
$ws.Range("I16").Value = 7000
$ws.Range("K16").Value = 7000
$ws.Range("M16").Value = -6713
# Row 29
$ws.Range("H29").Value = 16740
$ws.Range("I29").Value = 3000
$ws.Range("J29").Value = 20175
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 20175
$ws.Range("M29").Value = -2707
$ws.Range("N29").Value = -20761
# Row 31
$ws.Range("H31").Value = 13705.84
$ws.Range("I31").Value = 881.6875
$ws.Range("J31").Value = 36504.332
$ws.Range("K31").Value = 881.6875
$ws.Range("L31").Value = 36504.332
$ws.Range("M31").Value = -586.6875
$ws.Range("N31").Value = -37094.332
# Row 34
$ws.Range("H34").Value = 13705.84
$ws.Range("I34").Value = 881.6875
$ws.Range("J34").Value = 36504.332
$ws.Range("K34").Value = 881.6875
$ws.Range("L34").Value = 36504.332
$ws.Range("M34").Value = -679.6875
$ws.Range("N34").Value = -36908.332
# Row 113
$ws.Range("H113").Value = 11424.625
$ws.Range("I113").Value = 7000
$ws.Range("K113").Value = 7000
$ws.Range("M113").Value = -4830
# Row 134
$ws.Range("H134").Value = 35722164
$ws.Range("I134").Value = 3477.375
$ws.Range("K134").Value = 10432.125
$ws.Range("M134").Value = -7897.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 86
$ws.Range("H86").Value = 547.06665
$ws.Range("J86").Value = 829
$ws.Range("L86").Value = 2487
$ws.Range("N86").Value = -4859
# Row 89
$ws.Range("H89").Value = 547.06665
$ws.Range("J89").Value = 829
$ws.Range("L89").Value = 7461
$ws.Range("N89").Value = -19317

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2389.6924
$ws.Range("I113").Value = 1983.2727
$ws.Range("K113").Value = 1983.2727
$ws.Range("M113").Value = 186.7273
# Row 122
$ws.Range("H122").Value = 919020.4
$ws.Range("I122").Value = 1030150.2
$ws.Range("K122").Value = 3090450.6
$ws.Range("M122").Value = -3088000.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 45
$ws.Range("K2").Value = 45
$ws.Range("M2").Value = 67
# Row 16
$ws.Range("H16").Value = 62502856
$ws.Range("I16").Value = 90911840
$ws.Range("J16").Value = 3089.4
$ws.Range("K16").Value = 90911840
$ws.Range("L16").Value = 3089.4
$ws.Range("M16").Value = -90911670
$ws.Range("N16").Value = -3429.4
# Row 50
$ws.Range("H50").Value = 53893.332
$ws.Range("J50").Value = 53893.332
$ws.Range("L50").Value = 53893.332
$ws.Range("N50").Value = -55167.332
# Row 55
$ws.Range("H55").Value = 1921.3334
$ws.Range("I55").Value = 1503.6666
$ws.Range("K55").Value = 1503.6666
$ws.Range("M55").Value = -1330.6666
# Row 61
$ws.Range("H61").Value = 3263.923
$ws.Range("I61").Value = 3328.2
$ws.Range("J61").Value = 3049.6667
$ws.Range("K61").Value = 3328.2
$ws.Range("L61").Value = 3049.6667
$ws.Range("M61").Value = -3126.2
$ws.Range("N61").Value = -3453.6667
# Row 100
$ws.Range("H100").Value = 3453
$ws.Range("J100").Value = 4366.3335
$ws.Range("L100").Value = 4366.3335
$ws.Range("N100").Value = -5448.3335
# Row 113
$ws.Range("H113").Value = 3263.923
$ws.Range("I113").Value = 3328.2
$ws.Range("J113").Value = 3049.6667
$ws.Range("K113").Value = 3328.2
$ws.Range("L113").Value = 3049.6667
$ws.Range("M113").Value = -1158.2
$ws.Range("N113").Value = -7389.6667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 11264.706
$ws.Range("I62").Value = 8780
$ws.Range("K62").Value = 8780
$ws.Range("M62").Value = -8156
# Row 65
$ws.Range("H65").Value = 11264.706
$ws.Range("I65").Value = 8780
$ws.Range("K65").Value = 43900
$ws.Range("M65").Value = -40780
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 122
$ws.Range("H122").Value = 893170.4
$ws.Range("I122").Value = 1309809.4
$ws.Range("K122").Value = 3929428.2
$ws.Range("M122").Value = -3926978.2
# Row 136
$ws.Range("H136").Value = 16170.064
$ws.Range("I136").Value = 2508.5789
$ws.Range("K136").Value = 7525.736699999999
$ws.Range("M136").Value = -4975.736699999999
